# Generate Report for Archive
#
# Semantics of this edit (derived from the target diff):
#  - The handback status for "0a2f19d8-...md" and "ba847e91-...md" moved
#    from "Ready for handoff" to "In Translation" on every sheet
#    (Overview, zh-cn, de-de).
#  - As a consequence, the report re-lists the two "In Translation" rows
#    together: row 7 and row 8 swap their entire row content (all data
#    columns *and* the hyperlink display text), while rows 9 and 10 keep
#    their row position. Hyperlink relationship ids stay bound to the row
#    position (only TextToDisplay / cell values move).

$wb = $excel.ActiveWorkbook

function Set-CellAndLink {
    param(
        $ws,
        [string]$addr,
        [string]$value
    )
    $ws.Range($addr).Value = $value
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $value
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (zh-cn status), C (de-de
# status), D (Latest Handoff Date). Only A has hyperlinks.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndLink $wsOverview '$A$7' "0a2f19d8-b420-41fa-a9aa-761df885ee19.md"
$wsOverview.Range("B7").Value = "In Translation"
$wsOverview.Range("C7").Value = "In Translation"
$wsOverview.Range("D7").Value = "2016-03-24 12:33:27"

Set-CellAndLink $wsOverview '$A$8' "1dd44055-c2ac-4e8e-8401-a29cf035e1b0.md"
$wsOverview.Range("B8").Value = "In Translation"
$wsOverview.Range("C8").Value = "In Translation"
$wsOverview.Range("D8").Value = "2016-03-24 12:24:51"

$wsOverview.Range("B9").Value = "In Translation"
$wsOverview.Range("C9").Value = "In Translation"

# Row 10 (cd733652) is unchanged.

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A (Source File Name), C (Status), D (Latest
# Handoff File), E (Latest Handoff Datetime). A and D carry hyperlinks.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-CellAndLink $wsZh '$A$7' "0a2f19d8-b420-41fa-a9aa-761df885ee19.md"
$wsZh.Range("C7").Value = "In Translation"
Set-CellAndLink $wsZh '$D$7' "0a2f19d8-b420-41fa-a9aa-761df885ee19.7e9ffb89346cf519ff064dc99398873c91fe8639.zh-cn.xlf"
$wsZh.Range("E7").Value = "2016-03-24 12:33:22"

Set-CellAndLink $wsZh '$A$8' "1dd44055-c2ac-4e8e-8401-a29cf035e1b0.md"
$wsZh.Range("C8").Value = "In Translation"
Set-CellAndLink $wsZh '$D$8' "1dd44055-c2ac-4e8e-8401-a29cf035e1b0.18a4445dd0fa2b3cf7b50aa69ec47be7e66cceb6.zh-cn.xlf"
$wsZh.Range("E8").Value = "2016-03-24 12:24:19"

$wsZh.Range("C9").Value = "In Translation"

# Row 10 (cd733652) is unchanged.

# ---------------------------------------------------------------------
# Sheet "de-de": same layout as zh-cn.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-CellAndLink $wsDe '$A$7' "0a2f19d8-b420-41fa-a9aa-761df885ee19.md"
$wsDe.Range("C7").Value = "In Translation"
Set-CellAndLink $wsDe '$D$7' "0a2f19d8-b420-41fa-a9aa-761df885ee19.7e9ffb89346cf519ff064dc99398873c91fe8639.de-de.xlf"
$wsDe.Range("E7").Value = "2016-03-24 12:33:27"

Set-CellAndLink $wsDe '$A$8' "1dd44055-c2ac-4e8e-8401-a29cf035e1b0.md"
$wsDe.Range("C8").Value = "In Translation"
Set-CellAndLink $wsDe '$D$8' "1dd44055-c2ac-4e8e-8401-a29cf035e1b0.18a4445dd0fa2b3cf7b50aa69ec47be7e66cceb6.de-de.xlf"
$wsDe.Range("E8").Value = "2016-03-24 12:24:51"

$wsDe.Range("C9").Value = "In Translation"

# Row 10 (cd733652) is unchanged.
